$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 5).Value = "2026-02-18 06:19:32"
$ws.Cells.Item(2, 8).Formula = '="72%"'
$ws.Cells.Item(2, 8).Copy() | Out-Null
$ws.Cells.Item(2, 8).PasteSpecial(-4163) | Out-Null
$ws.Cells.Item(2, 14).Value = "-1.9 °C 5:35 TU"
$ws.Cells.Item(3, 5).Value = "2026-02-18 06:19:34"
$ws.Cells.Item(3, 8).Formula = '="96%"'
$ws.Cells.Item(3, 8).Copy() | Out-Null
$ws.Cells.Item(3, 8).PasteSpecial(-4163) | Out-Null
$ws.Cells.Item(4, 5).Value = "2026-02-18 06:19:37"
$ws.Cells.Item(4, 10).Value = "1018.4 hPa"
$ws.Cells.Item(4, 14).Value = "4.4 °C 5:50 TU"
$ws.Cells.Item(4, 15).Value = "7.1 °C"
$ws.Cells.Item(5, 5).Value = "2026-02-18 06:19:40"
$ws.Cells.Item(5, 8).Formula = '="88%"'
$ws.Cells.Item(5, 8).Copy() | Out-Null
$ws.Cells.Item(5, 8).PasteSpecial(-4163) | Out-Null
$ws.Cells.Item(5, 13).Value = "-0.1 °C 5:51 TU"
$ws.Cells.Item(5, 15).Value = "-2.0 °C"
$ws.Cells.Item(6, 5).Value = "2026-02-18 06:19:42"
$ws.Cells.Item(6, 10).Value = "1018.1 hPa"
$ws.Cells.Item(6, 15).Value = "8.0 °C"
$ws.Cells.Item(7, 5).Value = "2026-02-18 06:19:45"
$ws.Cells.Item(7, 10).Value = "1018.3 hPa"
$ws.Cells.Item(7, 14).Value = "11.4 °C 5:59 TU"
$ws.Cells.Item(7, 15).Value = "12.3 °C"
$ws.Cells.Item(8, 5).Value = "2026-02-18 06:19:48"
$ws.Cells.Item(8, 10).Value = "1018.4 hPa"
$ws.Cells.Item(8, 14).Value = "8.0 °C 5:30 TU"
$ws.Cells.Item(9, 5).Value = "2026-02-18 06:19:50"
$ws.Cells.Item(9, 14).Value = "4.0 °C 5:41 TU"
$ws.Cells.Item(9, 15).Value = "4.9 °C"
$ws.Cells.Item(10, 5).Value = "2026-02-18 06:19:53"
$ws.Cells.Item(10, 11).Value = "-0.1 MJ/m2"
$ws.Cells.Item(10, 12).Value = "5.8 km/h - 78º 5:35 TU"
$ws.Cells.Item(10, 15).Value = "6.7 °C"
$ws.Cells.Item(11, 5).Value = "2026-02-18 06:19:56"
$ws.Cells.Item(11, 14).Value = "-0.3 °C 5:44 TU"
$ws.Cells.Item(11, 15).Value = "1.9 °C"
$ws.Cells.Item(12, 5).Value = "2026-02-18 06:19:58"
$ws.Cells.Item(12, 14).Value = "4.4 °C 5:37 TU"
$ws.Cells.Item(12, 15).Value = "5.9 °C"
$ws.Cells.Item(13, 5).Value = "2026-02-18 06:20:01"
$ws.Cells.Item(13, 10).Value = "1023.1 hPa"
$ws.Cells.Item(13, 14).Value = "-3.8 °C 5:59 TU"
$ws.Cells.Item(13, 15).Value = "-2.3 °C"
$ws.Cells.Item(14, 5).Value = "2026-02-18 06:20:04"
$ws.Cells.Item(14, 12).Value = "11.9 km/h - 305º 5:53 TU"
$ws.Cells.Item(14, 15).Value = "9.9 °C"
$ws.Cells.Item(15, 5).Value = "2026-02-18 06:20:06"
$ws.Cells.Item(15, 14).Value = "4.3 °C 5:48 TU"
$ws.Cells.Item(16, 5).Value = "2026-02-18 06:20:09"
$ws.Cells.Item(16, 8).Formula = '="38%"'
$ws.Cells.Item(16, 8).Copy() | Out-Null
$ws.Cells.Item(16, 8).PasteSpecial(-4163) | Out-Null
$ws.Cells.Item(16, 15).Value = "0.6 °C"
$ws.Cells.Item(17, 5).Value = "2026-02-18 06:20:12"
$ws.Cells.Item(17, 8).Formula = '="88%"'
$ws.Cells.Item(17, 8).Copy() | Out-Null
$ws.Cells.Item(17, 8).PasteSpecial(-4163) | Out-Null
$ws.Cells.Item(17, 12).Value = "52.9 km/h - 280º 5:39 TU"
$ws.Cells.Item(17, 13).Value = "4.0 °C 5:40 TU"
$ws.Cells.Item(17, 15).Value = "2.2 °C"
$ws.Cells.Item(18, 5).Value = "2026-02-18 06:20:14"
$ws.Cells.Item(18, 10).Value = "1018.4 hPa"
$ws.Cells.Item(18, 14).Value = "4.6 °C 5:59 TU"
$ws.Cells.Item(18, 15).Value = "7.4 °C"
$ws.Cells.Item(19, 5).Value = "2026-02-18 06:20:17"
$ws.Cells.Item(19, 14).Value = "5.2 °C 5:31 TU"
$ws.Cells.Item(20, 5).Value = "2026-02-18 06:20:20"
$ws.Cells.Item(20, 8).Formula = '="78%"'
$ws.Cells.Item(20, 8).Copy() | Out-Null
$ws.Cells.Item(20, 8).PasteSpecial(-4163) | Out-Null
$ws.Cells.Item(20, 15).Value = "-1.1 °C"
$ws.Cells.Item(21, 5).Value = "2026-02-18 06:20:22"
$ws.Cells.Item(21, 10).Value = "1020.7 hPa"
$ws.Cells.Item(21, 14).Value = "0.4 °C 5:56 TU"
$ws.Cells.Item(21, 15).Value = "1.9 °C"
$ws.Cells.Item(22, 5).Value = "2026-02-18 06:20:25"
$ws.Cells.Item(23, 5).Value = "2026-02-18 06:20:27"
$ws.Cells.Item(23, 12).Value = "43.9 km/h - 307º 5:33 TU"
$ws.Cells.Item(23, 13).Value = "1.8 °C 5:45 TU"
$ws.Cells.Item(23, 15).Value = "0.5 °C"
$ws.Cells.Item(24, 5).Value = "2026-02-18 06:20:30"
$ws.Cells.Item(24, 8).Formula = '="100%"'
$ws.Cells.Item(24, 8).Copy() | Out-Null
$ws.Cells.Item(24, 8).PasteSpecial(-4163) | Out-Null
$ws.Cells.Item(24, 10).Value = "1018.9 hPa"
$ws.Cells.Item(24, 14).Value = "2.5 °C 5:59 TU"
$ws.Cells.Item(24, 15).Value = "5.1 °C"
$ws.Cells.Item(25, 5).Value = "2026-02-18 06:20:33"
$ws.Cells.Item(25, 13).Value = "1.4 °C 5:56 TU"
$ws.Cells.Item(25, 15).Value = "-0.3 °C"
$ws.Cells.Item(26, 5).Value = "2026-02-18 06:20:36"
$ws.Cells.Item(27, 5).Value = "2026-02-18 06:20:38"
$ws.Cells.Item(27, 8).Formula = '="51%"'
$ws.Cells.Item(27, 8).Copy() | Out-Null
$ws.Cells.Item(27, 8).PasteSpecial(-4163) | Out-Null
$ws.Cells.Item(28, 5).Value = "2026-02-18 06:20:41"
$ws.Cells.Item(28, 10).Value = "1018.8 hPa"
$ws.Cells.Item(28, 14).Value = "2.4 °C 5:55 TU"
$ws.Cells.Item(28, 15).Value = "4.7 °C"
$ws.Cells.Item(29, 5).Value = "2026-02-18 06:20:43"
$ws.Cells.Item(29, 8).Formula = '="92%"'
$ws.Cells.Item(29, 8).Copy() | Out-Null
$ws.Cells.Item(29, 8).PasteSpecial(-4163) | Out-Null
$ws.Cells.Item(29, 15).Value = "9.3 °C"
$ws.Cells.Item(30, 5).Value = "2026-02-18 06:20:46"
$ws.Cells.Item(30, 10).Value = "1018.4 hPa"
$ws.Cells.Item(30, 14).Value = "5.0 °C 5:44 TU"
$ws.Cells.Item(30, 15).Value = "5.9 °C"
$ws.Cells.Item(31, 5).Value = "2026-02-18 06:20:49"
$ws.Cells.Item(31, 8).Formula = '="78%"'
$ws.Cells.Item(31, 8).Copy() | Out-Null
$ws.Cells.Item(31, 8).PasteSpecial(-4163) | Out-Null
$ws.Cells.Item(31, 10).Value = "1017.0 hPa"
$ws.Cells.Item(31, 15).Value = "10.5 °C"
$ws.Cells.Item(32, 5).Value = "2026-02-18 06:20:51"
$ws.Cells.Item(32, 15).Value = "5.0 °C"
$ws.Cells.Item(33, 5).Value = "2026-02-18 06:20:54"
$ws.Cells.Item(33, 14).Value = "-1.5 °C 5:50 TU"
$ws.Cells.Item(33, 15).Value = "-0.1 °C"
$ws.Cells.Item(34, 5).Value = "2026-02-18 06:20:57"
$ws.Cells.Item(34, 8).Formula = '="53%"'
$ws.Cells.Item(34, 8).Copy() | Out-Null
$ws.Cells.Item(34, 8).PasteSpecial(-4163) | Out-Null
$ws.Cells.Item(34, 15).Value = "0.1 °C"
$ws.Cells.Item(35, 5).Value = "2026-02-18 06:20:59"
$ws.Cells.Item(35, 8).Formula = '="82%"'
$ws.Cells.Item(35, 8).Copy() | Out-Null
$ws.Cells.Item(35, 8).PasteSpecial(-4163) | Out-Null
$ws.Cells.Item(35, 14).Value = "4.7 °C 5:46 TU"
$ws.Cells.Item(35, 15).Value = "7.1 °C"
$ws.Cells.Item(36, 5).Value = "2026-02-18 06:21:02"
$ws.Cells.Item(37, 5).Value = "2026-02-18 06:21:05"
$ws.Cells.Item(37, 10).Value = "1021.2 hPa"
$ws.Cells.Item(37, 14).Value = "0.0 °C 5:52 TU"
$ws.Cells.Item(37, 15).Value = "1.3 °C"
$ws.Cells.Item(38, 5).Value = "2026-02-18 06:21:07"
$ws.Cells.Item(38, 14).Value = "5.9 °C 5:56 TU"
$ws.Cells.Item(38, 15).Value = "8.8 °C"
$ws.Cells.Item(39, 5).Value = "2026-02-18 06:21:10"
$ws.Cells.Item(39, 13).Value = "1.9 °C 5:35 TU"
$ws.Cells.Item(39, 15).Value = "0.3 °C"
$ws.Cells.Item(40, 5).Value = "2026-02-18 06:21:13"
$ws.Cells.Item(40, 10).Value = "1021.6 hPa"
$ws.Cells.Item(40, 14).Value = "-0.7 °C 5:58 TU"
$ws.Cells.Item(40, 15).Value = "0.7 °C"
$ws.Cells.Item(41, 5).Value = "2026-02-18 06:21:15"
$ws.Cells.Item(41, 10).Value = "1018.0 hPa"
$ws.Cells.Item(41, 11).Value = "-0.1 MJ/m2"
$ws.Cells.Item(42, 5).Value = "2026-02-18 06:21:18"
$ws.Cells.Item(42, 8).Formula = '="98%"'
$ws.Cells.Item(42, 8).Copy() | Out-Null
$ws.Cells.Item(42, 8).PasteSpecial(-4163) | Out-Null
$ws.Cells.Item(42, 14).Value = "5.3 °C 5:50 TU"
$ws.Cells.Item(42, 15).Value = "8.1 °C"
$ws.Cells.Item(43, 5).Value = "2026-02-18 06:21:21"
$ws.Cells.Item(43, 14).Value = "5.6 °C 5:42 TU"
$ws.Cells.Item(43, 15).Value = "7.0 °C"
$ws.Cells.Item(44, 5).Value = "2026-02-18 06:21:23"
$ws.Cells.Item(44, 8).Formula = '="65%"'
$ws.Cells.Item(44, 8).Copy() | Out-Null
$ws.Cells.Item(44, 8).PasteSpecial(-4163) | Out-Null
$ws.Cells.Item(44, 13).Value = "-1.1 °C 5:56 TU"
$ws.Cells.Item(44, 15).Value = "-3.5 °C"
$ws.Cells.Item(45, 5).Value = "2026-02-18 06:21:26"
$ws.Cells.Item(45, 7).Value = "2 cm"
$ws.Cells.Item(45, 10).Value = "1020.5 hPa"
$ws.Cells.Item(45, 14).Value = "-0.1 °C 5:32 TU"
$ws.Cells.Item(46, 5).Value = "2026-02-18 06:21:29"
$ws.Cells.Item(46, 10).Value = "1018.9 hPa"
$ws.Cells.Item(46, 15).Value = "6.2 °C"

$excel.CutCopyMode = 0

